$d = $word.ActiveDocument

# The page footer/boilerplate block that must be removed consists of three
# consecutive paragraphs right after the "Requisitos" section's course
# requirement line:
#   1. an empty paragraph
#   2. "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3. "(c) 2020 . Contact: luizeleno@usp.br. ..." copyright/footer line
#
# Locate paragraph 2 ("Ver no Jupiter ...") by its text, then expand the
# deletion range to cover the paragraph immediately before it (the blank
# separator paragraph) through the end of the paragraph immediately after
# it (the copyright line), and delete that whole range - paragraph marks
# included - in one shot.

$markerText = "Ver no Jupiter Salvar em pdf Salvar em docx"

$markerIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq $markerText) {
        $markerIndex = $i
        break
    }
}

if ($markerIndex -ge 0) {
    $firstToRemove = $d.Paragraphs.Item($markerIndex - 1)
    $lastToRemove = $d.Paragraphs.Item($markerIndex + 1)
    $removalRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
    $removalRange.Delete()
}

Write-Output "paragraphs.count=$($d.Paragraphs.Count)"
